$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update the CEIC add-in metadata blob stored in the A1 cell comment
$newCommentText = "rDMAAB+LCAAAAAAAAAPtW1tvG8cV/isLPjVAqN2lKJlSRxvwIqlEKMkQqdjKizHcHUlbL3fZnVlJfHMeAgOuXaCoYyBpm9YBigYBGvshLVw7aP9LYMrOU/9Cz1z2xotFum5gAwyMmHMuM2fOnMu3gzH64LznaackpG7gbxTMJaOgEd8OHNc/3ihE7KhorhY+sNDmuU28qzjEPcJAWAMtn66fU3ejcMJYf13Xz87Ols6Wl4LwWC8Zhqlf32m17RPSw0XXpwz7NikkWs7lWgUL1Z3eDmHYwQxLzY1Cs91cqhPXbgBtB/v4mIRLtYi6PqF002cucwnlmiHBjNQbOx/JjVmlpdUlE+lj9FSyFrmeI+VykpKu5GBZ0nF7xCoZZqVorBXNSsc01kvL6+XKUrm88nGsmAiiFqasTcJT1xaENsO9vlA31syKaZSWyxWkTxSCuVIHWGjPc/bJqUuJUyeeR+fyiK4OsGoz2PV8zjSQntFVE72+Cdsh7p90XOaR+czY36lpPf99rfXhPvytbEons9BWEBIb/Phapu2Ss71QubfTbwG3c+KGbNDAg7nnOqAk3OtzZ82naqFG4LOqR0J20IczJw6EBDAsFkYE6VOYqVLDpTb8dv2IONYR9mhWKcdE14LwJu1jm+xCPut8jjPfC7ADgcdcylw7XXSMga6GQR9mhMVrgedswaxKeAIjmbnpg4v5srUguJlaN4mJxKmK84Uz7WEWi4/RUfskONvzvUE76lI7dLvEadRi6Yk8xBNTadcjyoIeWJGSkKRlKAP4DxJxlIwaxHZ72LvqgROptQyz5AioGrHgyGX1wIt6Po1tGqGia7CjDjlPdpiM0R4crs+dHvhNP5aXbp7IyivsB2fJmuMM4YQMuUrt+LjHGaPCDaDFxzfOESfCd7nletAnsmeRoeajon1CCJsYEpKDeEnc4p3Hqg12o14X0qsLOXYqVqVIT/kI4hRiHeyyDOgmRfGnYxjr4g/YkbDRpu9Ml4uZCJbLrGWZwBshIdiTV/OwfxOo11x2sluN9zKBg6QHpsqP8xBkbt/DA0FOvJSloaZve5FDZEFo+kciRLlt8lCnstEYqQU5biHsDzqDPtRn6q4z+LFRgI69TlkImKBg2UHks3DAKwfSlehlOjTq+mIB7M2scxSSX0UARQZbkW/XA2f21RzpnQPfZbNbGEShLIezqwjv8coY0QbhNUYU/Zn17Xn2RMO5xHs+6QW+a8/ubXAyt955jY3QOKtm1iAyv2aW96Cvy7bHc31mtRBwJDS6uZapUhrYrghWlR5ORl+fkjINcoQjDzAcgxZ7nNTeUTKq0pujMlkSOgi9uAJaHCFTgMi201uyATxwGLhkBz1O0AGZXmsjPSvPEZBNNv3jFvaPI8AYSV0ZpSf1l/fHToh9yreTQIqRUjxZCMV1SkIdSxavvUgEgixeAXCRPiKHOqTXD0Ls7YBj3C0VdgovARTZwexEjaC3ecSOnaynqolW3rLY8MvERJOS2+AJr8rkCFEI8b1IMJ7KpDTEd7kDaenVsed2Q1lV41Y+iQcHloLDuP7yzc0JFOMzgA8y6L4fkgFH6elA0UXImjFDBjAvpFZ7v1wprRgrK4Bs+BiJHW8Hp0wDDEwAI2pFrYPPtZ9rTd9xT10nwh78hOCDnUn3qBY0l2ZWB23F9V2YUAU785S8AOCMYxfax7hgwkkVrEOCQ2+QEZRbbAU2yF3c//fw6R8u7v9j+Olfii++/s3wzsP/fP/F8yffPH/6VFLlHqU06uCuR4RFnVqlYiyXIcASEuJe1QUmdiKbCdrhoYDCyRipLzsxqG8269utmigkCTFWl71E5x+NgyBKh225C7GQOEs9DgEpYnXiwqTGOW6mN1n8I+6U5KWz/GmK0hcvnv31xbO/TdVWDktBlrm2VimapUsxGHwOm2NyCQZr5Yo/Fy4XjZViqZQRHpFB+7LyJ35qOtayaawZ5hXDTIq4k0TxJKFRlpqpg4/1ET1JqktclIRAdhwzReR3CGUJW+ZCZqBi9O+/fvnos5yU8q6i5GcB4wR+4Yvp8UBMvbvf0dp7B/v1Ta2z2eZxkvIycnLyVwir1ZOEygWV70OWv69BV4cuphXgU6igBUcawfaJNoBUzCRiLtgmUeVCrznlqJXbYRD15YlkFFLqBMmknEzUmFBsBE/4c6zqpKwJ4tLW4T+/m6SgNtJIkay6AEniUtBQjiNJGb7K2s//9fzJbV7cHv92+OST3AxqneQ+AOIcsik7TMIeSp5qNCMUdK0tnHnTuJFpLIrIv6CuBq7PqGWWxceTGiFQNfls4m/U7EGvExMLfwF9hIJ+genmOVOJbe0iPU8AO/sY2myQfnQmBFnDU7/++McvL37/3cWDxy9vfzO88/Xw7oMXz/708tuHMusuPnt8ce9bVeVHG4GwhX/KSvSniYsRW+PZqPGmrf1w63eaHzANsIYWiYr0w63PM5NxQwUqSWcGLJcYkjdhTDSrzPW0jCmJDTm9REV2/jpvYcuJhGpiQd+100U+LvKpeN4Jxs+anWJEiRYAjHoPdpIXTpVn1VMqsqVevWKUzJLiSmv4FrqYZly/7QVdgBExQ9w8jIjktF6tkMqK9bZbe7VqKxWRRuyFDgl5GMofKMaSvKU0aTyKQy1DAS4gPjvy+GXRmNg4K5k5U8Z0detyVHV4+Zt8R5GTQPUoDCUi8tVlfjvqAwqO7+am88V9ZQb47kqQmoXC6bjZyPNhnOFCI8yzOUHwRWlSLFmmmpRf8Egcu8tdkw6Bl7vjBHeo+3qJtABVklDndWczDINwYvFJObHYDkBoqCh66vFERpyphNtOelYxIS54bwp2XzGnwu7LwPXbCKEXYPlNgeUVY0awXPqJwHJlFrBcWYDlBVhegOWJYNmcDpbNG5l2MA6WV5cXYHkBlt8psGwuwPK7ApZHis/bDJZXc2D5I+xFRINzJI7m+vyOOaLQ0ifh5qmibweEfvnw7vMnXwy/+nJ458/DW98vUPSbQtGV8mwo2lj9aVB0aXUGFC2FFih6gaLfJhQdv7UbR9Ex55Uo+uLBox9v3x1++nh47+Hw3lcXj+7/j3C6NB1Ol25kGsY4nC4t4PQCTr9bcLq0gNPvCpweKT7/Fzitp49FkroYNIhH2HwPqfVUeyc4fW1dOPt5VZt0z3OUM+d775K4JZ0g+5qcB8qbekwug64ahoC0+OvTuV9/x6+m9rF/PKdVcktCkT8/gtXV8+wtN6TsOm9M6pekHCaUQwlZr/OHyPKHGB9aZVMSQEDPzq7nzIxTmMm3/YHXcnvunG+SjDjP85OAL/t9iema80UMbzG75BwQZ2YGKI7dX0L7kI/45plNBi7U1ESfP5yl7vEJm9ewK11MHNI1inaXlIplx6gU1whZLpom/B/bJfi2WOHPbtXkUEFccjbnInp8YOm/ubH+CyC6WfysMwAA"
$ws.Range("A1").Comment.Text($newCommentText)

# 3. Update the custom number format (numFmtId 166) used by B27:D42
$ws.Range("B27:D42").NumberFormat = "###0.000"

# 4. Relabel A11 from "Function Description" to "Function Information"
$ws.Range("A11").Value = "Function Information"

# 5. Adjust a handful of statistic cells to their slightly revised values
$ws.Range("B20").Value2 = 0.0156159028969747
$ws.Range("C20").Value2 = 0.5088159326089575
$ws.Range("B21").Value2 = -1.075751539480832
$ws.Range("C21").Value2 = -0.497756315366257
